$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.142.02"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.778.69"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.43"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.549"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.66"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0656"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "2.034.53"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("E13").Value = "  +6.33%  "
$ws.Range("D14").Value = "1.768.45"
$ws.Range("E14").Value = "  -3.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.625"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("D16").Value = "34.122.20"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.74"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.24"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("E23").Value = "  -3.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  -3.38%  "
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.40"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.99"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  -3.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0512"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").Value = "1.440.14"
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.623"
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.85"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.83"
$ws.Range("E40").Value = "  -3.39%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("E42").Value = "  -3.60%  "
$ws.Range("E43").Value = "  -5.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0510"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").Value = "1.934.90"
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.81"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.98"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.40"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.63"
$ws.Range("E51").Value = "  -6.64%  "
